$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "90.702.02"
$c.ClearFormats()
$ws.Range("E2").Value = "  +3.07%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.162.61"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  +0.17%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.04"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.81%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "625.05"
$c.ClearFormats()
$ws.Range("E6").Value = "  +2.01%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.392"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.19%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.702"
$c.ClearFormats()
$ws.Range("E8").Value = "  +4.10%  "

$ws.Range("E9").Value = "  +0.15%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.157.99"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.86%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.560"
$c.ClearFormats()
$ws.Range("E11").Value = "  +3.23%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.179"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.85%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000251"
$c.ClearFormats()
$ws.Range("E13").Value = "  +2.33%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "90.435.91"
$c.ClearFormats()
$ws.Range("E14").Value = "  +3.14%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.29"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.34%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.749.75"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.43%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "32.29"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.51%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.203.28"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.66%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.09%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000207"
$c.ClearFormats()
$ws.Range("E20").Value = "  +55.81%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "436.98"
$c.ClearFormats()
$ws.Range("E21").Value = "  +5.47%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "13.20"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.52%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.40"
$c.ClearFormats()
$ws.Range("E23").Value = "  -1.64%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.94"
$c.ClearFormats()
$ws.Range("E24").Value = "  -3.06%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.13"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.54%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.54"
$c.ClearFormats()
$ws.Range("E26").Value = "  -6.18%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "79.81"
$c.ClearFormats()
$ws.Range("E27").Value = "  +8.39%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "3.341.30"
$c.ClearFormats()
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  +0.27%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.154"
$c.ClearFormats()
$ws.Range("E31").Value = "  -4.99%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.ClearFormats()
$ws.Range("E32").Value = "  +31.16%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "8.21"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.81%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "518.37"
$c.ClearFormats()
$ws.Range("E34").Value = "  -5.65%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.86"
$c.ClearFormats()
$ws.Range("E35").Value = "  -2.26%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.87"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.03%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.27"
$c.ClearFormats()
$ws.Range("E37").Value = "  -4.76%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "22.20"
$c.ClearFormats()
$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("E39").Value = "  +2.49%  "

$ws.Range("E40").Value = "  +0.30%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.ClearFormats()
$ws.Range("E41").Value = "  -5.64%  "

$ws.Range("E42").Value = "  -0.02%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.90"
$c.ClearFormats()
$ws.Range("E43").Value = "  -2.13%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.365"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.73%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "145.37"
$c.ClearFormats()
$ws.Range("E45").Value = "  -1.85%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "44.01"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.57%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "169.56"
$c.ClearFormats()
$ws.Range("E47").Value = "  -3.03%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.ClearFormats()
$ws.Range("E48").Value = "  -1.36%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.736"
$c.ClearFormats()
$ws.Range("E49").Value = "  +5.53%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "24.35"
$c.ClearFormats()
$ws.Range("E50").Value = "  +1.51%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.20"
$c.ClearFormats()
$ws.Range("E51").Value = "  -3.52%  "
